$wb = $excel.ActiveWorkbook

$configWs = $wb.Worksheets.Item("config")

# Add the new "OutSensors" row that points at the external output-sensor list
$configWs.Range("A17").Value = "OutSensors"
$configWs.Range("A17").Font.Bold = $true
$configWs.Range("B17").Value = "./OutListParameters.xlsx"

# Make the config sheet the active sheet / tab, with A17 selected
$configWs.Activate()
$configWs.Range("A17").Select()
